# "Update on Michael Laptop"
#
# Student "Nick Miller" (row 11 of the schedule table) drops the
# "BUS 236 2 - Statistics for Business" course and picks up
# "PHI 223 - Intro to Formal Logic" instead (CRN 91602 -> 91528,
# meeting time 09:30-10:45A -> 12:00-12:50P). Only row 11 changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A11 (the CRN) is stored as text ("91528"), not a number, in the
# original workbook. Assigning the numeric-looking string directly
# would make Excel auto-convert it to a Number, so instead build it as
# a text formula and flatten it back down to a literal with a
# copy/paste-values round trip -- that keeps the cell's text type (and
# its existing style) intact, just like the source file.
$ws.Range("A11").Formula = "=""91528"""
$ws.Range("A11").Copy()
$ws.Range("A11").PasteSpecial(-4163)

$ws.Range("B11").Value = "PHI 223"
$ws.Range("C11").Value = "Intro to Formal Logic"
$ws.Range("E11").Value = "12:00-12:50P"
$ws.Range("F11").Value = "12:00-12:50P"
